$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '22.445.65', '  -0.01%  '),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.569.98', '  +0.24%  '),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.002', '  +0.04%  '),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  +0.07%  '),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '288.58', '  +0.01%  '),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.3704', '  +1.00%  '),
    @('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '48.31', '  -3.33%  '),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3308', '  -1.84%  '),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.135', '  -0.19%  '),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07498', '  +0.31%  '),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  +0.04%  '),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.71', '  -1.34%  '),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.927', '  -0.95%  '),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.870', '  -1.09%  '),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.569.26', '  +0.46%  '),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001118', '  +0.67%  '),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06748', '  +0.38%  '),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '87.47', '  -2.80%  '),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.001', '  +0.00%  '),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.343', '  -0.39%  '),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '16.59', '  +2.47%  '),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '12.02', '  +0.01%  '),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '22.444.22', '  +0.05%  '),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.386', '  -0.19%  '),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.587', '  -1.24%  '),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '153.67', '  +2.69%  '),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.69', '  -0.53%  '),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '5.024', '  -0.08%  '),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '124.34', '  +0.27%  '),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.745.62', '  +0.44%  '),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.055', '  +0.32%  '),
    @('WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.014', '  +0.33%  '),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.112', '  -0.48%  '),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.788', '  +1.93%  '),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08338', '  +0.63%  '),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02467', '  +1.19%  '),
    @('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2259', '  -0.07%  '),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06402', '  -0.05%  '),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.285', '  -3.88%  '),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.333', '  -0.13%  '),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.32', '  +1.43%  '),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6317', '  +2.80%  '),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '13.81', '  -0.47%  '),
    @('Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6151', '  +6.68%  '),
    @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.770', '  +0.11%  '),
    @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.060', '  +1.41%  '),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '126.20', '  +0.63%  '),
    @('EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.213', '  -0.77%  '),
    @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.07219', '  -1.59%  '),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '76.85', '  +2.41%  '),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[$i][2]
    $dCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
